# Update the ObjTables/SBtab header metadata (version bump 0.0.9 -> 1.0.0,
# date 2020-04-27 01:08:22/23 -> 2020-05-29 00:22:17) stamped into the first
# one or two header rows of every SBtab sheet in the workbook.
$wb = $excel.ActiveWorkbook

# First sheet ("!!Compartment") carries the workbook-wide "!!!ObjTables" banner
# in A1 in addition to its own "!!ObjTables ... class='Compartment'" line in A2.
$ws = $wb.Worksheets.Item("!!Compartment")
$ws.Unprotect()
$ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='1.0.0' date='2020-05-29 00:22:17'"
$ws.Range("A2").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Compartment' name='Compartment' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Compound")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Compound' name='Compound' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Definition")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Definition' name='Definition' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Enzyme")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Enzyme' name='Enzyme' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!FbcObjective")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='FbcObjective' name='FbcObjective' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Gene")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Gene' name='Gene' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Layout")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Layout' name='Layout' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Measurement")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Measurement' name='Measurement' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!PbConfig")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='PbConfig' name='PbConfig' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Position")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Position' name='Position' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Protein")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Protein' name='Protein' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Quantity")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Quantity' name='Quantity' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!QuantityInfo")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='QuantityInfo' name='QuantityInfo' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!QuantityMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='QuantityMatrix' name='QuantityMatrix' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Reaction")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Reaction' name='Reaction' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!ReactionStoichiometry")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Regulator")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Regulator' name='Regulator' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Relation")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Relation' name='Relation' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Relationship")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='Relationship' name='Relationship' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrix' name='SparseMatrix' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixColumn")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixOrdered")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixRow")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='SparseMatrixRow' name='SparseMatrixRow' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!StoichiometricMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!rxnconContingencyList")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='rxnconContingencyList' name='rxnconContingencyList' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!rxnconReactionList")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='rxnconReactionList' name='rxnconReactionList' date='2020-05-29 00:22:17' objTablesVersion='1.0.0'"
$ws.Protect()
